$d = $word.ActiveDocument

# --- 1. Insert a new "Assertions" heading paragraph right after the
#        paragraph that contains the last inline image (do this before
#        touching NoProofing below, so the freshly split run does not
#        inherit <w:noProof/> from its image-bearing neighbour).
$shapeCount = $d.InlineShapes.Count
$lastShape = $d.InlineShapes.Item($shapeCount)
$imgParagraph = $lastShape.Range.Paragraphs.Item(1)

# Remember its position in the Paragraphs collection so we can grab the
# freshly-created paragraph right after inserting it.
$imgParaStart = $imgParagraph.Range.Start
$imgParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Start -eq $imgParaStart) {
        $imgParaIndex = $i
        break
    }
}

$imgParagraph.Range.InsertParagraphAfter()

$newParagraph = $d.Paragraphs.Item($imgParaIndex + 1)
$newRange = $newParagraph.Range
$newRange.InsertAfter("Assertions")
$newRange.Font.Size = 28
$newRange.Font.SizeBi = 28

# --- 2. Mark the three screenshot drawings (InlineShapes 2-4) as "no
#        proofing", i.e. add <w:rPr><w:noProof/></w:rPr> to the run that
#        hosts each <w:drawing>. InlineShape 1 already carries noProof in
#        the source document, so we only touch the later three images.
for ($i = 2; $i -le $shapeCount; $i++) {
    $shp = $d.InlineShapes.Item($i)
    $shp.Range.NoProofing = -1
}
